# Fix lỗi trong report cơ sở — update "Lương" sheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 428571.4285714286
$ws.Range("B13").Value = 285714.2857142857
$ws.Range("B23").Value = 428571.4285714286
$ws.Range("B31").Value = 428571.4285714286
$ws.Range("B32").Value = 285714.2857142857
$ws.Range("B33").Value = 428571.4285714286
$ws.Range("B34").Value = 1142857.142857143
